$d = $word.ActiveDocument

# --- Paragraph 1: title ---------------------------------------------------
# "Still Time To Protest The Draft" / Heading1  ->  same text / Title style
$p1 = $d.Paragraphs.Item(1)
$p1.Style = "Title"

# --- Paragraph 2: author line ----------------------------------------------
# "By Dorothy Day" (bold, no style) -> "Dorothy Day" (Authors style, not bold)
#
# Rather than editing the existing run in place (which leaves a stray,
# now-redundant <w:rPr> behind once the bold direct-formatting is cleared),
# insert a brand-new paragraph right after the title, style it first, and
# then type the author text into it - that way the run never carries the
# old bold formatting at all. The original "By Dorothy Day" paragraph is
# then removed.
$p1.Range.Collapse(0)                 # wdCollapseEnd
$p1.Range.InsertParagraphAfter()

$authors = $d.Paragraphs.Item(2)
$authors.Style = "Authors"
$authors.Range.Text = "Dorothy Day"

$oldAuthors = $d.Paragraphs.Item(3)
$oldRange = $oldAuthors.Range
$oldRange.Expand(4)                   # wdParagraph - grab the mark too
$oldRange.Delete()
